# "Added more to projects"
# Insert two new bulleted list items (same list/style as the existing
# "Ideas to begin with:" list, numId 3) right after the "App that provides
# recipes..." item and before the trailing empty bookmark paragraph.

$d = $word.ActiveDocument

# Locate the "App that provides recipes..." list paragraph to anchor on.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "App that provides recipes*") {
        $anchor = $p
        break
    }
}

# First new bullet: condense-English app idea.
$anchor.Range.InsertParagraphAfter()
$firstNew = $anchor.Next()
$firstNew.Range.Text = "App to condense English in documents such as reports, emails or job applications."

# Second new bullet: "Baby brain" (curly quotes).
$firstNew.Range.InsertParagraphAfter()
$secondNew = $firstNew.Next()
$secondNew.Range.Text = [char]0x201C + "Baby brain" + [char]0x201D
